$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2062.5715
$ws.Range("I62").Value = 1799.6
$ws.Range("K62").Value = 1799.6
$ws.Range("M62").Value = -1175.6
$ws.Range("H65").Value = 2062.5715
$ws.Range("I65").Value = 1799.6
$ws.Range("K65").Value = 8998
$ws.Range("M65").Value = -5878
$ws.Range("H70").Value = 125887.5
$ws.Range("I70").Value = 1016.6667
$ws.Range("K70").Value = 3050.0001
$ws.Range("M70").Value = -2780.0001
$ws.Range("H73").Value = 125887.5
$ws.Range("I73").Value = 1016.6667
$ws.Range("K73").Value = 3050.0001
$ws.Range("M73").Value = -2114.0001
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120
$ws.Range("H137").Value = 916.05
$ws.Range("I137").Value = 989.5294
$ws.Range("K137").Value = 2968.5882
$ws.Range("M137").Value = -418.5882000000001
$ws.Range("H138").Value = 4540.263
$ws.Range("I138").Value = 2857.9412
$ws.Range("J138").Value = 5025
$ws.Range("K138").Value = 8573.8236
$ws.Range("L138").Value = 15075
$ws.Range("M138").Value = -3433.8236
$ws.Range("N138").Value = -25355

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4934.6885
$ws.Range("I32").Value = 477.57352
$ws.Range("J32").Value = 38610.668
$ws.Range("K32").Value = 477.57352
$ws.Range("L32").Value = 38610.668
$ws.Range("M32").Value = -190.57352
$ws.Range("N32").Value = -39184.668
$ws.Range("H45").Value = 14825.75
$ws.Range("I45").Value = 20017.908
$ws.Range("K45").Value = 20017.908
$ws.Range("M45").Value = -19640.908
$ws.Range("H61").Value = 4569.778
$ws.Range("I61").Value = 4569.778
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4569.778
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4357.778
$ws.Range("N61").ClearContents()
$ws.Range("H122").Value = 2617.5833
$ws.Range("I122").Value = 2447.1
$ws.Range("K122").Value = 7341.299999999999
$ws.Range("M122").Value = -4891.299999999999
$ws.Range("H132").Value = 4729.5713
$ws.Range("I132").Value = 4766.577
$ws.Range("J132").Value = 4248.5
$ws.Range("K132").Value = 14299.731
$ws.Range("L132").Value = 12745.5
$ws.Range("M132").Value = -11769.731
$ws.Range("N132").Value = -17805.5
$ws.Range("H136").Value = 4569.778
$ws.Range("I136").Value = 4569.778
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13709.334
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -11159.334
$ws.Range("N136").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 42860656
$ws.Range("I86").Value = 50003600
$ws.Range("K86").Value = 50003600
$ws.Range("M86").Value = -50002477
$ws.Range("H89").Value = 42860656
$ws.Range("I89").Value = 50003600
$ws.Range("K89").Value = 250018000
$ws.Range("M89").Value = -250012384
$ws.Range("H107").Value = 3329.975
$ws.Range("I107").Value = 1774.9474
$ws.Range("J107").Value = 4736.905
$ws.Range("K107").Value = 1774.9474
$ws.Range("L107").Value = 4736.905
$ws.Range("M107").Value = 145.0526
$ws.Range("N107").Value = -8576.904999999999
$ws.Range("H134").Value = 3182.652
$ws.Range("I134").Value = 2533.9048
$ws.Range("K134").Value = 7601.714399999999
$ws.Range("M134").Value = -5066.714399999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3083.4216
$ws.Range("I31").Value = 2058.5789
$ws.Range("K31").Value = 2058.5789
$ws.Range("M31").Value = -1763.5789
$ws.Range("H34").Value = 3083.4216
$ws.Range("I34").Value = 2058.5789
$ws.Range("K34").Value = 2058.5789
$ws.Range("M34").Value = -1856.5789
$ws.Range("H44").Value = 10000
$ws.Range("I44").Value = 10000
$ws.Range("K44").Value = 10000
$ws.Range("M44").Value = -9558
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9155
$ws.Range("H58").Value = 5630.875
$ws.Range("I58").Value = 5630.875
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 5630.875
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5427.875
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value = 6862.6875
$ws.Range("I62").Value = 5964.1816
$ws.Range("J62").Value = 8839.4
$ws.Range("K62").Value = 5964.1816
$ws.Range("L62").Value = 8839.4
$ws.Range("M62").Value = -5340.1816
$ws.Range("N62").Value = -10087.4
$ws.Range("H65").Value = 6862.6875
$ws.Range("I65").Value = 5964.1816
$ws.Range("J65").Value = 8839.4
$ws.Range("K65").Value = 29820.908
$ws.Range("L65").Value = 44197
$ws.Range("M65").Value = -26700.908
$ws.Range("N65").Value = -50437
$ws.Range("H107").Value = 1701.9565
$ws.Range("I107").Value = 1584.8572
$ws.Range("J107").Value = 1884.1111
$ws.Range("K107").Value = 1584.8572
$ws.Range("L107").Value = 1884.1111
$ws.Range("M107").Value = 335.1428000000001
$ws.Range("N107").Value = -5724.1111
$ws.Range("H136").Value = 5630.875
$ws.Range("I136").Value = 5630.875
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 16892.625
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -14342.625
$ws.Range("N136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 619.4783
$ws.Range("I23").Value = 527.5714
$ws.Range("J23").Value = 659.6875
$ws.Range("K23").Value = 1582.7142
$ws.Range("L23").Value = 1979.0625
$ws.Range("M23").Value = -1347.7142
$ws.Range("N23").Value = -2449.0625
$ws.Range("H68").Value = 1799.75
$ws.Range("I68").Value = 1626.5714
$ws.Range("J68").Value = 1893
$ws.Range("K68").Value = 4879.7142
$ws.Range("L68").Value = 5679
$ws.Range("M68").Value = -4068.7142
$ws.Range("N68").Value = -7301
$ws.Range("H71").Value = 1799.75
$ws.Range("I71").Value = 1626.5714
$ws.Range("J71").Value = 1893
$ws.Range("K71").Value = 14639.1426
$ws.Range("L71").Value = 17037
$ws.Range("M71").Value = -10583.1426
$ws.Range("N71").Value = -25149
$ws.Range("H132").Value = 3921.0688
$ws.Range("J132").Value = 3944.4783
$ws.Range("L132").Value = 35500.3047
$ws.Range("N132").Value = -40560.3047

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4999.067
$ws.Range("I70").Value = 4900.8
$ws.Range("K70").Value = 4900.8
$ws.Range("M70").Value = -4630.8
$ws.Range("H73").Value = 4999.067
$ws.Range("I73").Value = 4900.8
$ws.Range("K73").Value = 4900.8
$ws.Range("M73").Value = -3964.8
$ws.Range("H122").Value = 5929.8423
$ws.Range("I122").Value = 4124.4546
$ws.Range("J122").Value = 8412.25
$ws.Range("K122").Value = 12373.3638
$ws.Range("L122").Value = 25236.75
$ws.Range("M122").Value = -9923.363799999999
$ws.Range("N122").Value = -30136.75
$ws.Range("H132").Value = 2967.238
$ws.Range("I132").Value = 2967.238
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8901.714
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6371.714
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9272.125
$ws.Range("I40").Value = 6804.3335
$ws.Range("J40").Value = 13385.111
$ws.Range("K40").Value = 6804.3335
$ws.Range("L40").Value = 13385.111
$ws.Range("M40").Value = -6668.3335
$ws.Range("N40").Value = -13657.111
$ws.Range("H61").Value = 2775
$ws.Range("J61").Value = 1995.4286
$ws.Range("L61").Value = 1995.4286
$ws.Range("N61").Value = -2399.4286
$ws.Range("H82").Value = 1236.0769
$ws.Range("I82").Value = 910.2857
$ws.Range("J82").Value = 1616.1666
$ws.Range("K82").Value = 910.2857
$ws.Range("L82").Value = 1616.1666
$ws.Range("M82").Value = -549.2857
$ws.Range("N82").Value = -2338.1666
$ws.Range("H85").Value = 1236.0769
$ws.Range("I85").Value = 910.2857
$ws.Range("J85").Value = 1616.1666
$ws.Range("K85").Value = 910.2857
$ws.Range("L85").Value = 1616.1666
$ws.Range("M85").Value = 337.7143
$ws.Range("N85").Value = -4112.1666
$ws.Range("H113").Value = 2775
$ws.Range("J113").Value = 1995.4286
$ws.Range("L113").Value = 1995.4286
$ws.Range("N113").Value = -6335.4286
$ws.Range("H132").Value = 23035.406
$ws.Range("I132").Value = 32236.52
$ws.Range("K132").Value = 96709.56
$ws.Range("M132").Value = -94179.56

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 500
$ws.Range("I113").Value = 500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670
